# docs: writing pipeline documentation
#
# Update the subtitle on slide 1: the conference location is clarified
# from "l'Ecole Centrale" to "l'Ecole Centrale (Paris)".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$target = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        if ($shp.TextFrame.TextRange.Text -like "*Ecole Centrale*") {
            $target = $shp
            break
        }
    }
}

if ($target -eq $null) {
    # Fallback: the subtitle placeholder is the second shape on the slide.
    $target = $s.Shapes.Item(2)
}

$target.TextFrame.TextRange.Text = "Conférence donnée à l’Ecole Centrale (Paris) en novembre 2022"
